# Add two new combat-dialog lines after "YOU HUNTED MY KIND FOR YEARS! NOW ITS MY TURN!"
# Both new lines use the same ListParagraph/bulleted-list formatting as the
# surrounding dialog entries.

$d = $word.ActiveDocument

# Locate the paragraph whose text is the existing last dialog line.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "YOU HUNTED MY KIND FOR YEARS! NOW ITS MY TURN!") {
        $target = $p
    }
}

# Insert a new (empty) paragraph right after it; InsertParagraphAfter clones
# the paragraph formatting (pStyle + numPr) of $target, matching the diff.
$target.Range.InsertParagraphAfter()
$newPara1 = $target.Next()
$newPara1.Range.Text = "YOU THINK ITS HARD TO KILL ME NOW YOU SHOULD HAVE SEEN ME WITH ALL OF MY LIMBS!!!!!!!!!"

# Insert a second new paragraph after the first new one, same treatment.
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "WE HAVE BEEN FIGHTING SO LONG MY LIMBS ARE GETTING RUSTY!"
